# Apply weekly report updates for WR_89708709_WeekEnding_062925.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates ---

# Report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"

# Total Billed Amount
$ws.Range("C8").Value = 2219.58

# Total Line Items
$ws.Range("C9").Value = 4

# Scope ID # is now blank
$ws.Range("G10").Value = ""

# --- Friday (06/27/2025) section pricing ---
$ws.Range("H16").Value = 858.75
$ws.Range("H17").Value = 631.14
$ws.Range("H18").Value = 648.53
$ws.Range("H19").Value = 2138.42

# --- Saturday (06/28/2025) section ---
# Row 24 used to be "CNA-TR" / "CNA,Transfer Conductor" with 5 units and 0 pricing.
# It is replaced with the former row 25's line item (CNC-SNB-40), updated with
# the real billed amount, and the old row 25 is removed entirely (rows shift up).
$ws.Range("B24").Value = "CNC-SNB-40"
$ws.Range("D24").Value = "CNC,Splice Non-Tension Bare,#1/0-#4/0"
$ws.Range("F24").Value = 4
$ws.Range("H24").Value = 81.16

# Remove the now-duplicate line item row (old row 25); this shifts the TOTAL
# row (old row 26) up to row 25 and fixes the A:G merge + dimension/mergeCells
# bookkeeping automatically.
$ws.Rows("25:25").Delete()

# The TOTAL for the Saturday section (now row 25) reflects the single
# remaining line item.
$ws.Range("H25").Value = 81.16
